# Add I0 and IF columns (I and J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1), same style as other header cells (e.g. H1 -> style index 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for rows 2..74 -> columns I (col 9) and J (col 10)
$data = @(
  @(2, 8, 8),
  @(3, 9, 9),
  @(4, 9, 9),
  @(5, 9, 9),
  @(6, 8, 9),
  @(7, 9, 9),
  @(8, 9, 9),
  @(9, 9, 9),
  @(10, 10, 11),
  @(11, 9, 9),
  @(12, 9, 9),
  @(13, 9, 9),
  @(14, 9, 10),
  @(15, 9, 9),
  @(16, 8, 9),
  @(17, 9, 10),
  @(18, 9, 9),
  @(19, 9, 9),
  @(20, 9, 9),
  @(21, 10, 10),
  @(22, 9, 10),
  @(23, 9, 9),
  @(24, 7, 7),
  @(25, 9, 9),
  @(26, 9, 9),
  @(27, 9, 9),
  @(28, 9, 9),
  @(29, 9, 9),
  @(30, 9, 9),
  @(31, 9, 9),
  @(32, 9, 10),
  @(33, 9, 9),
  @(34, 9, 9),
  @(35, 9, 9),
  @(36, 9, 9),
  @(37, 9, 9),
  @(38, 9, 9),
  @(39, 9, 9),
  @(40, 9, 9),
  @(41, 9, 9),
  @(42, 8, 8),
  @(43, 9, 9),
  @(44, 9, 9),
  @(45, 9, 9),
  @(46, 9, 9),
  @(47, 9, 9),
  @(48, 9, 9),
  @(49, 9, 9),
  @(50, 9, 9),
  @(51, 9, 9),
  @(52, 9, 9),
  @(53, 10, 10),
  @(54, 9, 9),
  @(55, 9, 9),
  @(56, 9, 9),
  @(57, 9, 9),
  @(58, 9, 9),
  @(59, 9, 9),
  @(60, 10, 10),
  @(61, 9, 9),
  @(62, 9, 9),
  @(63, 9, 9),
  @(64, 9, 9),
  @(65, 9, 9),
  @(66, 9, 9),
  @(67, 6, 6),
  @(68, 5, 5),
  @(69, 1, 1),
  @(70, 6, 6),
  @(71, 5, 5),
  @(72, 6, 6),
  @(73, 5, 5),
  @(74, 4, 4)
)

foreach ($row in $data) {
  $r = $row[0]
  $i0 = $row[1]
  $if = $row[2]
  $ws.Cells.Item($r, 9).Value = $i0
  $ws.Cells.Item($r, 10).Value = $if
}
